# Add new rows (A:1-15, B:6-20, C:<reaction numbers>) to both worksheets,
# reusing the existing formatted-number style from row 2 column A.

$wb = $excel.ActiveWorkbook

$nbrC = @(879, 866, 865, 921, 923, 917, 917, 915, 914, 0, 0, 896, 896, 901, 885)
$barC = @(1037, 1031, 1019, 940, 941, 930, 923, 917, 913, 0, 0, 912, 886, 886, 885)

$ws1 = $wb.Worksheets.Item(1)
for ($i = 0; $i -lt $nbrC.Length; $i++) {
    $r = 3 + $i
    $ws1.Range("A2").Copy()
    $ws1.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws1.Cells.Item($r, 1).Value = $i + 1
    $ws1.Cells.Item($r, 2).Value = $i + 6
    $ws1.Cells.Item($r, 3).Value = $nbrC[$i]
}

$ws2 = $wb.Worksheets.Item(2)
for ($i = 0; $i -lt $barC.Length; $i++) {
    $r = 3 + $i
    $ws2.Range("A2").Copy()
    $ws2.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws2.Cells.Item($r, 1).Value = $i + 1
    $ws2.Cells.Item($r, 2).Value = $i + 6
    $ws2.Cells.Item($r, 3).Value = $barC[$i]
}
